# Append latest scrape run (2026-02-05 06:55:59) to the "ランサーズ" sheet.
# New run returned fewer (5) listings than the previous run (19), so the
# sheet is rewritten with the new top-5 rows and truncated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-02-05 06:55:59"

# Final data for rows 2..6 (columns A-H), in order.
$rows = @(
    @{ B = "Amazonから情報取得するツール作ってください。SP-API有 Python希望";
       D = "10,000 円 ~ 20,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5486242";
       G = 435;
       H = "🔥Python,API ◆ツール" },
    @{ B = "SaaSビジネスにおける「バーティカル(垂直型)」展開の横スライド可能なAIシステムの開発です";
       D = "100,000 円 ~ 200,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5485911";
       G = 383;
       H = "🔥AI,Ai ◆開発" },
    @{ B = "【急募】クリックポスト自動発行ツール開発依頼";
       D = "20,000 円 ~ 50,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5485895";
       G = 123;
       H = "◆ツール,開発" },
    @{ B = "【Excelでのマクロ作成】リサーチツールの作成【スクレイピング】";
       D = "1,000 ~ 5,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5486225";
       G = 100;
       H = "◆ツール,スクレイピング" },
    @{ B = "【急募】iOS/AndroidアプリのSkyWay切替対応エンジニア募集";
       D = "100,000 円 ~ 200,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5486110";
       G = 38;
       H = "◇アプリ" }
)

# Category and due-date are constant across all rows (unchanged in diff).
$category = "システム開発"
$deadline = "期限情報なし"

# 1) Fix up the Hyperlinks collection FIRST, while rows 2-20 still sit at
#    their original positions (deleting rows later does not renumber or
#    drop the Hyperlinks collection, so this must happen before the row
#    delete below). Rows 2-6 get their target address updated in place
#    (keeps the existing Hyperlink cell style); rows 7-20 lose their
#    hyperlink registration entirely since those rows are being removed.
for ($r = 2; $r -le 6; $r++) {
    $newUrl = $rows[$r - 2].F
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq ("`$F`$" + $r)) {
            $h.Address = $newUrl
        }
    }
}

# NOTE: the Hyperlinks collection re-indexes after every single .Delete(),
# so collecting several objects first and deleting them in a batch causes
# some to be skipped. Deleting strictly one-at-a-time against a freshly
# re-fetched live collection (re-querying $ws.Hyperlinks.GetEnumerator on
# every pass) is the only pattern that reliably empties the target set.
$keepAddrs = @()
for ($r = 2; $r -le 6; $r++) { $keepAddrs += ("`$F`$" + $r) }

$changed = $true
while ($changed) {
    $changed = $false
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if (-not ($keepAddrs -contains $addr)) {
            $h.Delete()
            $changed = $true
            break
        }
    }
}

# 2) Delete the old extra rows (7-20), from the bottom up.
for ($r = 20; $r -ge 7; $r--) {
    $ws.Rows.Item($r).Delete()
}

# 3) Overwrite rows 2-6 with the new data.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $category
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $deadline
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
}

# 4) Column width adjustments (autofit-style tweaks that came along with
#    the new data).
$ws.Columns.Item(2).ColumnWidth = 48.17
$ws.Columns.Item(4).ColumnWidth = 27.17
$ws.Columns.Item(8).ColumnWidth = 17.17
